# Daily attendance processing - 2025-10-08 16:49:34
# Updates recorder-email orderings, attendance fractions, and computed
# percentage metrics on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Year 3 / C1 / ANATOMY session 2) ---
$ws.Range("G3").Value = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 4 (Year 3 / C1 / ANATOMY session 3) ---
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H4").Value = "133/221"

# --- Row 10 (Year 3 / C1 summary: Average Attendance %) ---
$ws.Range("L10").Value = "'39.5%"

# --- Row 12 (Year 3 / C1 / HISTOLOGY session 1) ---
$ws.Range("G12").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 15 (Group Statistics: Year 3 / C1) ---
$ws.Range("S15").Value = "'46.0%"

# --- Row 16 (Group Statistics: Year 3 / C2) ---
$ws.Range("S16").Value = "'34.6%"

# --- Row 25 (Year 3 / C2 / ANATOMY session 2) ---
$ws.Range("G25").Value = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 26 (Year 3 / C2 / ANATOMY session 3) ---
$ws.Range("G26").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H26").Value = "106/246"

# --- Row 34 (Year 3 / C2 / HISTOLOGY session 1) ---
$ws.Range("G34").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 41 (Year 3 / C2 / PHYSIOLOGY session 1) ---
$ws.Range("G41").Value = "neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
